$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (match style of existing header row, e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data values for columns I (I0) and J (IF)
$values = @(
    @(8, 9),
    @(6, 9),
    @(6, 8),
    @(9, 9),
    @(5, 8),
    @(5, 5),
    @(1, 2)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
